$p = $ppt.ActivePresentation

# Locate the "Заключение" slide (slide 7) whose body placeholder holds the
# paragraph about the BunBricky project, then switch that paragraph's font
# to Lucida Console (matching the other re-fonted text already present on
# this deck, e.g. the slide title and the paragraph's own trailing runs).
$targetSlide = $null
$targetShape = $null

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.HasTextFrame) {
            $text = $shape.TextFrame.TextRange.Text
            if ($text.Contains("BunBricky") -and $text.Contains("реализованы")) {
                $targetSlide = $slide
                $targetShape = $shape
            }
        }
    }
}

$textRange = $targetShape.TextFrame.TextRange
$paragraph = $textRange.Paragraphs(1, 1)
$paragraph.Font.Name = "Lucida Console"
